# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Row -> new value for the "展览" sheet (sheet1)
$exhibitUpdates = @{
    2  = 248
    3  = 1091
    5  = 432
    7  = 559
    8  = 71
    9  = 6808
    11 = 97
    12 = 142
    15 = 1100
    16 = 16217
    17 = 1593
    19 = 330
    21 = 116
    22 = 11364
    24 = 1007
    25 = 4478
    26 = 318
    28 = 47
    31 = 141
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value for the "全部类型" sheet (sheet4)
$allTypesUpdates = @{
    2  = 248
    3  = 1091
    5  = 432
    7  = 559
    9  = 71
    10 = 6808
    12 = 97
    13 = 142
    17 = 1100
    18 = 16217
    19 = 1593
    21 = 330
    23 = 116
    26 = 11364
    28 = 1007
    29 = 4478
    30 = 318
    32 = 47
    35 = 141
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
